$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "02_03_24" (first sheet) — sprint log updates
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("02_03_24")
$ws1.Activate()

# Row 14 ("FPS Counter Class") gains a Notes entry and becomes taller
# (wraps to two lines) now that the FPS counter was switched to the
# singleton pattern but isn't fully wired up yet.
$ws1.Range("E14").Value = "Displays currently, but does not update in real-time."
$ws1.Rows("14:14").RowHeight = 30

# Row 10 ("Setup Project") gains a Status + Notes entry.
$ws1.Range("D10").Value = "Implemented"
$ws1.Range("E10").Value = "One Solution."

# ---------------------------------------------------------------------
# Sheet "Remaining Tasks" — clarified notes header
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Remaining Tasks")
$ws3.Activate()
$ws3.Range("A1").Value = "Note: Each point is ~2 hours of work. Status can be 'Started,' 'Implemented,' 'Cleanup' and 'Finished.'"
$ws3.Range("E9").Select()

# ---------------------------------------------------------------------
# Back to "02_03_24" — add the new "Median Filter (3 Color)" task row.
# ---------------------------------------------------------------------
$ws1.Activate()

# Insert a brand-new data row (new task) directly below row 14, followed
# by a blank spacer row, before the existing "Ends/Total" summary rows.
$ws1.Rows("15:17").Insert()

# Copy the date-formatted style from row 14 onto the new row 15 so the
# date cell renders the same way as the rows above it.
$ws1.Range("A14").Copy()
$ws1.Range("A15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws1.Range("A15").Value = $ws1.Range("A14").Value2
$ws1.Range("B15").Value = 3
$ws1.Range("C15").Value = "Median Filter (3 Color)"
$ws1.Range("D15:E15").Clear()

# Rows 16-17 are a blank gap between the task table and the summary rows.
$ws1.Range("A16:E17").Clear()

$ws1.Range("D19").Select()

# Restore original active sheet/tab.
$ws1.Activate()
